$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Round 1 results: fill in column E (winner) for each pairing in rows 19-33.
# For most pairs the winner is the team already shown in column D (top of pair);
# for the 27/28 pairing, the winner is the team from the bottom row (D28, "Quokka").
$ws.Range("E19").Value = $ws.Range("D19").Value2
$ws.Range("E21").Value = $ws.Range("D21").Value2
$ws.Range("E23").Value = $ws.Range("D23").Value2
$ws.Range("E25").Value = $ws.Range("D25").Value2
$ws.Range("E27").Value = $ws.Range("D28").Value2
$ws.Range("E29").Value = $ws.Range("D29").Value2
$ws.Range("E31").Value = $ws.Range("D31").Value2
$ws.Range("E33").Value = $ws.Range("D33").Value2

# Update view state to match the author's final scroll/selection position.
$ws.Range("E33").Select()
$excel.ActiveWindow.ScrollRow = 16
